# Coinranking "cryptos" snapshot refresh (GitHub Actions symbol-list update).
#
# Most rows just get a refreshed Price (column D). A handful of rows also had
# their Coin / Link / Volume(1h) text (columns B, C, E) replaced because the
# coin that occupied that ranking slot changed between scrapes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds the "Price" as text (e.g. "247.44"), not a number. Excels
# COM layer auto-converts a numeric-looking string typed into a General cell
# into a real number, so mark these cells as Text first to preserve the
# original string formatting (leading/trailing zeros, etc.) exactly.
$priceCells = @(
    "D2", "D3", "D4", "D5", "D6", "D7", "D8", "D11", "D12", "D13", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D48", "D50"
)
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Price updates --------------------------------------------------------
$ws.Range('D2').Value = '247.64'  # was '247.44'
$ws.Range('D3').Value = '22.77'  # was '22.79'
$ws.Range('D4').Value = '5.556'  # was '5.572'
$ws.Range('D5').Value = '0.05642'  # was '0.05638'
$ws.Range('D6').Value = '3.404'  # was '3.403'
$ws.Range('D7').Value = '6.482'  # was '6.483'
$ws.Range('D8').Value = '0.8014'  # was '0.8017'
$ws.Range('D11').Value = '0.07413'  # was '0.07409'
$ws.Range('D12').Value = '0.03197'  # was '0.03178'
$ws.Range('D13').Value = '0.02966'  # was '0.02975'
$ws.Range('D15').Value = '0.001671'  # was '0.001669'
$ws.Range('D16').Value = '0.04707'  # was '2.964'
$ws.Range('D17').Value = '0.006270'  # was '0.04712'
$ws.Range('D18').Value = '0.001049'  # was '0.0005744'
$ws.Range('D19').Value = '0.003822'  # was '0.006267'
$ws.Range('D20').Value = '0.0001501'  # was '0.001053'
$ws.Range('D21').Value = '0.0004602'  # was '0.003822'
$ws.Range('D22').Value = '3.982'  # was '0.0001501'
$ws.Range('D23').Value = '2.084'  # was '0.0004603'
$ws.Range('D24').Value = '0.01173'  # was '3.981'
$ws.Range('D25').Value = '0.3311'  # was '2.112'
$ws.Range('D26').Value = '0.1291'  # was '0.3311'
$ws.Range('D27').Value = '2.074'  # was '0.1277'
$ws.Range('D40').Value = '0.04183'  # was '0.04181'
$ws.Range('D41').Value = '0.007143'  # was '0.007133'
$ws.Range('D42').Value = '0.003501'  # was '0.1046'
$ws.Range('D43').Value = '0.1045'  # was '0.002971'
$ws.Range('D44').Value = '0.008682'  # was '0.008698'
$ws.Range('D45').Value = '0.00005649'  # was '0.00005641'
$ws.Range('D47').Value = '0.6803'  # was '0.6805'
$ws.Range('D48').Value = '0.02817'  # was '0.02808'
$ws.Range('D50').Value = '0.01010'  # was '0.01011'

# --- Row content swaps (Coin / Link / Volume label) -----------------------
# Row 16
$ws.Range('B16').Value = 'CoinExToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('E16').Value = '15CoinExTokenCET'
# Row 17
$ws.Range('B17').Value = 'TigerCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('E17').Value = '16TigerCashTCH'
# Row 18
$ws.Range('B18').Value = 'BitKan'
$ws.Range('C18').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('E18').Value = '17BitKanKAN'
# Row 19
$ws.Range('B19').Value = 'HotbitToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('E19').Value = '18HotbitTokenHTB'
# Row 20
$ws.Range('B20').Value = 'NitroEx'
$ws.Range('C20').Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range('E20').Value = '19NitroExNTX'
# Row 21
$ws.Range('B21').Value = 'UpBots'
$ws.Range('C21').Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range('E21').Value = '20UpBotsUBXT'
# Row 22
$ws.Range('B22').Value = 'LEO'
$ws.Range('C22').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('E22').Value = '21LEOLEO'
# Row 23
$ws.Range('B23').Value = 'BTSEToken'
$ws.Range('C23').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('E23').Value = '22BTSETokenBTSE'
# Row 24
$ws.Range('B24').Value = 'One'
$ws.Range('C24').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('E24').Value = '23OneONEBestin24h'
# Row 25
$ws.Range('B25').Value = 'BitpandaEcosystemToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('E25').Value = '24BitpandaEcosystemTokenBEST'
# Row 26
$ws.Range('B26').Value = 'ProBitToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range('E26').Value = '25ProBitTokenPROB'
# Row 27
$ws.Range('B27').Value = 'MCDex'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('E27').Value = '26MCDexMCBWorstin24h'
# Row 42
$ws.Range('B42').Value = 'CEJI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('E42').Value = '41CEJICEJI'
# Row 43
$ws.Range('B43').Value = 'BKEXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('E43').Value = '42BKEXTokenBKK'
# Row 48
$ws.Range('E48').Value = '47BOLOBOLO'
